# Update the Maltaspor fantasy-basketball roster sheet: the player table
# (columns A:C, rows 2-19) is re-sorted into a new order, "Keon Ellis" is
# dropped from the roster, and "Dalano Banton" is added in his place - all
# other players keep their existing Position/Team pairing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("Derrick White","PG,SG","Boston Celtics"),
  @("Luke Kennard","SG","Memphis Grizzlies"),
  @("Damian Lillard","PG","Milwaukee Bucks"),
  @("Cade Cunningham","PG,SG","Detroit Pistons"),
  @("Devin Vassell","SG,SF","San Antonio Spurs"),
  @("Onyeka Okongwu","PF,C","Atlanta Hawks"),
  @("Malik Monk","PG,SG,SF","Sacramento Kings"),
  @("Naz Reid","PF,C","Minnesota Timberwolves"),
  @("Anthony Davis","PF,C","Los Angeles Lakers"),
  @("Julius Randle","PF,C","Minnesota Timberwolves"),
  @("Collin Sexton","PG,SG","Utah Jazz"),
  @("LaMelo Ball","PG,SG","Charlotte Hornets"),
  @("Dalano Banton","SG,SF","Portland Trail Blazers"),
  @("Deandre Ayton","C","Portland Trail Blazers"),
  @("Alex Caruso","SG,SF","Oklahoma City Thunder"),
  @("Coby White","PG,SG","Chicago Bulls"),
  @("Isaiah Hartenstein","C","Oklahoma City Thunder"),
  @("Cameron Johnson","SF,PF","Brooklyn Nets")
)

$r = 2
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $r = $r + 1
}
